$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5562
$ws1.Range("F6").Value = 76
$ws1.Range("F8").Value = 912
$ws1.Range("F10").Value = 2492
$ws1.Range("F12").Value = 112
$ws1.Range("F15").Value = 9
$ws1.Range("F16").Value = 2329
$ws1.Range("F17").Value = 282

# Sheet "全部类型" (sheet4): same underlying events, different row offsets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5562
$ws4.Range("F7").Value = 76
$ws4.Range("F10").Value = 912
$ws4.Range("F12").Value = 2492
$ws4.Range("F14").Value = 112
$ws4.Range("F18").Value = 9
$ws4.Range("F19").Value = 2329
$ws4.Range("F20").Value = 282
